$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''65.015.19'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +1.01%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''3.419.51'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +2.89%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  +0.17%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''545.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +2.87%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''178.15'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +0.63%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''0.631'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  +6.46%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '''  +0.05%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.620'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  +1.44%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''0.149'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  +8.45%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''53.01'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -1.75%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = '''  +2.99%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''9.09'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  +0.85%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''3.981.36'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  +3.17%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = '''  +2.10%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''3.428.94'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +3.35%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''18.13'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  +3.24%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''65.190.25'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +1.07%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''11.71'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +3.83%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''0.974'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  +1.08%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''412.06'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  +7.18%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''3.97'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  +6.65%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = '''  +1.91%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''84.25'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +2.44%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = '''  -4.38%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''2.83'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  +2.97%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''6.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -1.47%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''12.04'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  +5.91%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''8.81'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +6.39%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''29.54'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  +1.66%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''609.82'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -5.44%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '''  -4.12%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''11.55'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  +2.48%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = '''  +1.72%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''58.49'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  +1.69%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = '''Dai'
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = '''https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = '''0.998'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -0.17%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = '''Kaspa'
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = '''https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = '''0.147'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  +15.55%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''36.93'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +0.63%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''0.0₃0770'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +2.07%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''0.374'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -2.45%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''3.145.80'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  +5.03%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''3.31'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  +1.89%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '''  +0.19%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = '''  -4.41%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''2.77'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  +3.10%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''2.71'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +1.27%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''0.0406'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +0.78%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = '''  +1.42%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''  +3.27%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''137.99'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -0.38%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''8.31'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  +0.38%  '
$ws.Range("E51").Style = "Normal"
